# Update crypto price/volume figures per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.162.36"
$ws.Cells.Item(2, 5).Value = "  +0.45%  "

$ws.Cells.Item(3, 4).Value = "1.833.57"
$ws.Cells.Item(3, 5).Value = "  +0.08%  "

$ws.Cells.Item(4, 4).Value = "'0.9989"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

$ws.Cells.Item(5, 4).Value = "'241.74"
$ws.Cells.Item(5, 5).Value = "  +0.94%  "

$ws.Cells.Item(6, 4).Value = "'0.6571"
$ws.Cells.Item(6, 5).Value = "  -1.98%  "

$ws.Cells.Item(7, 4).Value = "'0.9997"
$ws.Cells.Item(7, 5).Value = "  -0.05%  "

$ws.Cells.Item(8, 4).Value = "'0.07398"
$ws.Cells.Item(8, 5).Value = "  -0.13%  "

$ws.Cells.Item(9, 4).Value = "'0.2921"
$ws.Cells.Item(9, 5).Value = "  -1.02%  "

$ws.Cells.Item(10, 4).Value = "'22.90"
$ws.Cells.Item(10, 5).Value = "  +0.83%  "

$ws.Cells.Item(11, 4).Value = "'0.07739"
$ws.Cells.Item(11, 5).Value = "  +1.27%  "

$ws.Cells.Item(12, 4).Value = "1.833.90"
$ws.Cells.Item(12, 5).Value = "  +0.04%  "

$ws.Cells.Item(13, 4).Value = "'4.986"
$ws.Cells.Item(13, 5).Value = "  -0.33%  "

$ws.Cells.Item(14, 4).Value = "'0.6650"
$ws.Cells.Item(14, 5).Value = "  -1.07%  "

$ws.Cells.Item(15, 4).Value = "'82.75"
$ws.Cells.Item(15, 5).Value = "  -3.51%  "

$ws.Cells.Item(16, 5).Value = "  -0.31%  "

$ws.Cells.Item(17, 4).Value = "'0.000008420"
$ws.Cells.Item(17, 5).Value = "  +2.31%  "

$ws.Cells.Item(18, 4).Value = "29.167.26"
$ws.Cells.Item(18, 5).Value = "  +0.45%  "

$ws.Cells.Item(19, 4).Value = "2.102.82"
$ws.Cells.Item(19, 5).Value = "  +1.03%  "

$ws.Cells.Item(20, 4).Value = "'226.90"
$ws.Cells.Item(20, 5).Value = "  -0.20%  "

$ws.Cells.Item(21, 4).Value = "'12.45"
$ws.Cells.Item(21, 5).Value = "  +0.22%  "

$ws.Cells.Item(22, 5).Value = "  +0.12%  "

$ws.Cells.Item(23, 4).Value = "'7.127"
$ws.Cells.Item(23, 5).Value = "  -2.41%  "

$ws.Cells.Item(24, 4).Value = "'0.9997"
$ws.Cells.Item(24, 5).Value = "  -0.04%  "

$ws.Cells.Item(25, 4).Value = "'158.52"
$ws.Cells.Item(25, 5).Value = "  -1.20%  "

$ws.Cells.Item(26, 4).Value = "'8.592"
$ws.Cells.Item(26, 5).Value = "  -0.81%  "

$ws.Cells.Item(27, 4).Value = "'0.1390"
$ws.Cells.Item(27, 5).Value = "  -2.84%  "

$ws.Cells.Item(28, 5).Value = "  -0.12%  "

$ws.Cells.Item(29, 5).Value = "  +1.03%  "

$ws.Cells.Item(30, 4).Value = "'4.113"
$ws.Cells.Item(30, 5).Value = "  -2.62%  "

$ws.Cells.Item(31, 4).Value = "'4.042"
$ws.Cells.Item(31, 5).Value = "  -1.70%  "

$ws.Cells.Item(32, 4).Value = "'1.192"
$ws.Cells.Item(32, 5).Value = "  -0.29%  "

$ws.Cells.Item(33, 4).Value = "'0.05258"
$ws.Cells.Item(33, 5).Value = "  -2.00%  "

$ws.Cells.Item(34, 4).Value = "'1.865"
$ws.Cells.Item(34, 5).Value = "  +0.80%  "

$ws.Cells.Item(35, 4).Value = "'0.7360"
$ws.Cells.Item(35, 5).Value = "  -2.04%  "

$ws.Cells.Item(36, 4).Value = "'1.142"
$ws.Cells.Item(36, 5).Value = "  +1.65%  "

$ws.Cells.Item(37, 4).Value = "'2.652"
$ws.Cells.Item(37, 5).Value = "  -1.08%  "

$ws.Cells.Item(38, 4).Value = "1.301.34"
$ws.Cells.Item(38, 5).Value = "  +0.95%  "

$ws.Cells.Item(39, 5).Value = "  -0.90%  "

$ws.Cells.Item(40, 4).Value = "'2.729"
$ws.Cells.Item(40, 5).Value = "  +0.77%  "

$ws.Cells.Item(41, 4).Value = "'0.9192"
$ws.Cells.Item(41, 5).Value = "  -0.25%  "

$ws.Cells.Item(42, 4).Value = "'0.08682"
$ws.Cells.Item(42, 5).Value = "  +8.41%  "

$ws.Cells.Item(43, 4).Value = "'5.960"
$ws.Cells.Item(43, 5).Value = "  -0.73%  "

$ws.Cells.Item(44, 4).Value = "'0.9993"
$ws.Cells.Item(44, 5).Value = "  +0.12%  "

$ws.Cells.Item(45, 4).Value = "'102.31"
$ws.Cells.Item(45, 5).Value = "  -1.87%  "

$ws.Cells.Item(46, 4).Value = "1.997.04"
$ws.Cells.Item(46, 5).Value = "  +0.96%  "

$ws.Cells.Item(47, 4).Value = "'0.5137"
$ws.Cells.Item(47, 5).Value = "  -0.76%  "

$ws.Cells.Item(48, 5).Value = "  -1.13%  "

$ws.Cells.Item(49, 4).Value = "'1.748"
$ws.Cells.Item(49, 5).Value = "  -0.13%  "

$ws.Cells.Item(50, 4).Value = "'63.21"
$ws.Cells.Item(50, 5).Value = "  -0.38%  "

$ws.Cells.Item(51, 4).Value = "'0.05845"
$ws.Cells.Item(51, 5).Value = "  -1.25%  "
